# Catchment Area Service Mapping - September: "Add files via upload" edit
#
# Summary of changes applied:
#  1. "Other" sheet: remove the sparse "No" rows (rows where the
#     "other services" flag is No and no detail columns are populated),
#     keeping only the fully populated "Yes" rows. Excel shifts the
#     remaining rows up and renumbers them, shrinking the sheet from
#     A1:K20 to A1:K10. Delete bottom-up so row numbers of
#     not-yet-processed blocks stay valid.
#  2. "Other" sheet: selection moves to C6.
#  3. "food_cc" sheet: column A is widened (manual resize, so the
#     bestFit flag is dropped) to ~37.5 (closest reachable width to the
#     recorded 37.54296875).
#  4. "health_cc" sheet: selection moves to B3.
#
# The "Other" sheet must remain the active/selected sheet at the end,
# so sheets that only need a selection/width change are touched first
# and "Other" is activated last.

$wb = $excel.ActiveWorkbook

# --- food_cc: widen column A, drop AutoFit/bestFit sizing ---------------
$foodCc = $wb.Worksheets.Item("food_cc")
$foodCc.Columns.Item(1).ColumnWidth = 36.6666666

# --- health_cc: move the selection to B3 ---------------------------------
$healthCc = $wb.Worksheets.Item("health_cc")
$healthCc.Activate()
$healthCc.Range("B3").Select()

# --- Other: delete the sparse "No" rows (bottom-up) -----------------------
$other = $wb.Worksheets.Item("Other")
$other.Activate()

$other.Rows("18:18").Delete()
$other.Rows("14:15").Delete()
$other.Rows("9:11").Delete()
$other.Rows("3:6").Delete()

# --- Other: finally move the selection to C6 and keep it the active sheet -
$other.Range("C6").Select()
